$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row at 94 (pushes the existing "合計" row down to 95)
# and populate it with the new daily data point (2020-04-28).
$ws.Rows.Item(94).Insert()
$ws.Range("A94").Value = 43949
$ws.Range("B94").Value = 396
$ws.Range("C94").Value = 31114
$ws.Range("D94").Value = 148
$ws.Range("E94").Value = 6664

# Update the selected cell shown in the saved view.
$ws.Range("E96").Select()

# Extend the print area to cover the newly added row.
$wb.Names.Item(1).RefersTo = "=相談件数!`$A`$1:`$E`$97"
